$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOTAS")

# New row 27: next subject, "DISEÑO DE SISTEMAS"
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 2031
$ws.Range("C27").Value = "DISEÑO DE SISTEMAS"
$ws.Range("D27").Value = 44265
$ws.Range("E27").Value = 10

# Extend the average formula to include the new row
$ws.Range("F2").Formula = "=(SUM(E2:E27))/A26"

# Mirror the formatting already used on row 26 (number format, borders, etc.)
$ws.Range("A26:E26").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Extend conditional formatting to cover the new row
$fcs = $ws.Range("E2:E26").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("E2:E27"))
}

# Update the selection like the author's last click before saving
$ws.Range("B29").Select()

$wb.Save()
